{"js": "// Task 4 asn Task 5: Adding a test suite\n//\n// In the paragraph describing testGetPrice_NegativeNumber(), insert the\n// parenthetical \"(I placed it under comments) \" right after\n// \"does not occur \" and before \"because the \", turning:\n//   \"...does not occur because the VendingMachineItem trhows...\"\n// into:\n//   \"...does not occur (I placed it under comments) because the VendingMachineItem trhows...\"\n// The automatic \"_GoBack\" bookmark (Word's last-edit marker) is moved to sit\n// right after the newly typed \"comments\" text, matching where the cursor\n// would be after typing the insertion.\n\nconst body = context.document.body;\n\n// 1. Locate the still-intact \"because the\" phrase and insert the new text\n//    immediately before it.\nconst hits = body.search(\"because the\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error(\"Anchor phrase 'because the' not found\");\n}\n\nhits.items[0].insertText(\"(I placed it under comments) \", \"Before\");\nawait context.sync();\n\n// 2. Move the \"_GoBack\" bookmark so it sits right after the newly inserted\n//    \"comments\" text (i.e. right before the closing parenthesis), matching\n//    where Word leaves it after typing.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst markHits = body.search(\"comments\", { matchCase: true });\nmarkHits.load(\"items\");\nawait context.sync();\n\nif (markHits.items.length === 0) {\n  throw new Error(\"Inserted text 'comments' not found\");\n}\n\nconst caret = markHits.items[0].getRange(\"After\");\ncaret.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Task 4 asn Task 5: Adding a test suite\n#\n# In the paragraph describing testGetPrice_NegativeNumber(), insert the\n# parenthetical \"(I placed it under comments) \" right after\n# \"does not occur \" and before \"because the \", turning:\n#   \"...does not occur because the VendingMachineItem trhows...\"\n# into:\n#   \"...does not occur (I placed it under comments) because the VendingMachineItem trhows...\"\n# The automatic \"_GoBack\" bookmark (Word's last-edit marker) is moved to sit\n# right after the newly typed \"comments\" text, matching where the cursor\n# would be after typing the insertion.\n\n$d = $word.ActiveDocument\n\n# 1. Locate the still-intact \"because the\" phrase and insert the new text\n#    immediately before it.\n$find = $d.Content\n$find.Find.Text = \"because the\"\n$found = $find.Find.Execute()\nif (-not $found) {\n    throw \"Anchor phrase 'because the' not found\"\n}\n$find.InsertBefore(\"(I placed it under comments) \")\n\n# 2. Move the \"_GoBack\" bookmark so it sits right after the newly inserted\n#    \"comments\" text (i.e. right before the closing parenthesis), matching\n#    where Word leaves it after typing.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$markRange = $d.Content\n$markRange.Find.Text = \"comments\"\n$markFound = $markRange.Find.Execute()\nif (-not $markFound) {\n    throw \"Inserted text 'comments' not found\"\n}\n$caret = $d.Range($markRange.End, $markRange.End)\n$d.Bookmarks.Add(\"_GoBack\", $caret)\n"}
